$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old emoji values (📗 / 📘) with the new ones (✅ / ⚠️)
$ws.Range("A2").Value = "✅"
$ws.Range("A3").Value = "⚠️"
$ws.Range("A4").Value = "⚠️"
